$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.767.14"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "2.445.08"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.32%  "
$ws.Range("E15").Value = "  +6.30%  "
$ws.Range("D16").Value = "2.889.90"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "62.686.88"
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "2.439.94"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  +9.42%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("B26").Value = "Bittensor"
$ws.Range("C26").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "650.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.17%  "
$ws.Range("B27").Value = "BabyDogeCoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D27").Value = "0.0₆0602"
$ws.Range("E27").Value = "  +116.16%  "
$ws.Range("E28").Value = "  +17.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").Value = "  +6.09%  "
$ws.Range("D31").Value = "2.571.60"
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.40%  "
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.62%  "
$ws.Range("E44").Value = "  +5.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("E50").Value = "  +7.21%  "
$ws.Range("E51").Value = "  +2.78%  "
